$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 25014692
$ws.Range("J113").Value = 33352068
$ws.Range("L113").Value = 33352068
$ws.Range("N113").Value = -33358576

$ws.Range("H132").Value = 1459.0652
$ws.Range("I132").Value = 1264.881
$ws.Range("J132").Value = 3498
$ws.Range("K132").Value = 3794.643
$ws.Range("L132").Value = 10494
$ws.Range("M132").Value = -1264.643
$ws.Range("N132").Value = -15554

$ws.Range("H137").Value = 8664.828
$ws.Range("J137").Value = 11019.786
$ws.Range("L137").Value = 33059.358
$ws.Range("N137").Value = -38159.358

$ws.Range("H138").Value = 2910.3132
$ws.Range("I138").Value = 2212
$ws.Range("J138").Value = 3098.3206
$ws.Range("K138").Value = 6636
$ws.Range("L138").Value = 9294.961800000001
$ws.Range("M138").Value = -1496
$ws.Range("N138").Value = -19574.9618

$ws.Range("H141").Value = 1882.625
$ws.Range("I141").Value = 1882.625
$ws.Range("K141").Value = 5647.875
$ws.Range("M141").Value = -467.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2005583.2
$ws.Range("I32").Value = 2025336.6
$ws.Range("K32").Value = 2025336.6
$ws.Range("M32").Value = -2025049.6

$ws.Range("H61").Value = 32262368
$ws.Range("I61").Value = 2732.2104
$ws.Range("J61").Value = 83340130
$ws.Range("K61").Value = 2732.2104
$ws.Range("L61").Value = 83340130
$ws.Range("M61").Value = -2520.2104
$ws.Range("N61").Value = -83340554

$ws.Range("H74").Value = 4585.353
$ws.Range("I74").Value = 3240.2
$ws.Range("K74").Value = 3240.2
$ws.Range("M74").Value = -2366.2

$ws.Range("H77").Value = 4585.353
$ws.Range("I77").Value = 3240.2
$ws.Range("K77").Value = 16201
$ws.Range("M77").Value = -11833

$ws.Range("H110").Value = 18525220
$ws.Range("I110").Value = 7476.5625
$ws.Range("K110").Value = 7476.5625
$ws.Range("M110").Value = -5431.5625

$ws.Range("H122").Value = 4711.143
$ws.Range("I122").Value = 3706.8572
$ws.Range("J122").Value = 5715.4287
$ws.Range("K122").Value = 11120.5716
$ws.Range("L122").Value = 17146.2861
$ws.Range("M122").Value = -8670.571599999999
$ws.Range("N122").Value = -22046.2861

$ws.Range("H132").Value = 3703.6597
$ws.Range("I132").Value = 2412.7812
$ws.Range("K132").Value = 7238.3436
$ws.Range("M132").Value = -4708.3436

$ws.Range("H136").Value = 32262368
$ws.Range("I136").Value = 2732.2104
$ws.Range("J136").Value = 83340130
$ws.Range("K136").Value = 8196.6312
$ws.Range("L136").Value = 250020390
$ws.Range("M136").Value = -5646.6312
$ws.Range("N136").Value = -250025490

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9811063
$ws.Range("I20").Value = 15156734
$ws.Range("J20").Value = 10664.667
$ws.Range("K20").Value = 15156734
$ws.Range("L20").Value = 10664.667
$ws.Range("M20").Value = -15156487
$ws.Range("N20").Value = -11158.667

$ws.Range("H107").Value = 56252188
$ws.Range("I107").Value = 59212620
$ws.Range("J107").Value = 4001
$ws.Range("K107").Value = 59212620
$ws.Range("L107").Value = 4001
$ws.Range("M107").Value = -59210700
$ws.Range("N107").Value = -7841

$ws.Range("H113").Value = 5068.4287
$ws.Range("I113").Value = 5068.4287
$ws.Range("K113").Value = 5068.4287
$ws.Range("M113").Value = -2898.4287

$ws.Range("H128").Value = 3610.1
$ws.Range("I128").Value = 3610.1
$ws.Range("K128").Value = 10830.3
$ws.Range("M128").Value = -8340.299999999999

$ws.Range("H134").Value = 3626223.2
$ws.Range("I134").Value = 4718810.5
$ws.Range("K134").Value = 14156431.5
$ws.Range("M134").Value = -14153896.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H76").Value = 5341.875
$ws.Range("I76").Value = 5341.875
$ws.Range("K76").Value = 5341.875
$ws.Range("M76").Value = -5026.875

$ws.Range("H79").Value = 5341.875
$ws.Range("I79").Value = 5341.875
$ws.Range("K79").Value = 5341.875
$ws.Range("M79").Value = -4249.875

$ws.Range("H94").Value = 1441.5555
$ws.Range("I94").Value = 1725.125
$ws.Range("J94").Value = 1214.7
$ws.Range("K94").Value = 1725.125
$ws.Range("L94").Value = 1214.7
$ws.Range("M94").Value = -1274.125
$ws.Range("N94").Value = -2116.7

$ws.Range("H132").Value = 5159.925
$ws.Range("I132").Value = 3730.9614
$ws.Range("J132").Value = 7813.7144
$ws.Range("K132").Value = 11192.8842
$ws.Range("L132").Value = 23441.1432
$ws.Range("M132").Value = -8662.8842
$ws.Range("N132").Value = -28501.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 17157788
$ws.Range("I4").Value = 19112634
$ws.Range("J4").Value = 1239758.4
$ws.Range("K4").Value = 57337902
$ws.Range("L4").Value = 3719275.2
$ws.Range("M4").Value = -57337790
$ws.Range("N4").Value = -3719499.2

$ws.Range("H26").Value = 682.3333
$ws.Range("I26").Value = 189.42857
$ws.Range("J26").Value = 1372.4
$ws.Range("K26").Value = 568.28571
$ws.Range("L26").Value = 4117.200000000001
$ws.Range("M26").Value = -280.28571
$ws.Range("N26").Value = -4693.200000000001

$ws.Range("H114").Value = 1295.2307
$ws.Range("I114").Value = 314.75
$ws.Range("J114").Value = 1731
$ws.Range("K114").Value = 944.25
$ws.Range("L114").Value = 5193
$ws.Range("M114").Value = 2309.75
$ws.Range("N114").Value = -11701

$ws.Range("H117").Value = 2500
$ws.Range("J117").Value = 2500
$ws.Range("L117").Value = 7500
$ws.Range("N117").Value = -14384

$ws.Range("H128").Value = 145618.5
$ws.Range("I128").Value = 145618.5
$ws.Range("K128").Value = 436855.5
$ws.Range("M128").Value = -431875.5

$ws.Range("H131").Value = 1815.7709
$ws.Range("I131").Value = 1485.8572
$ws.Range("J131").Value = 1872.0975
$ws.Range("K131").Value = 4457.571599999999
$ws.Range("L131").Value = 5616.2925
$ws.Range("M131").Value = 582.4284000000007
$ws.Range("N131").Value = -15696.2925

$ws.Range("H137").Value = 5052.846
$ws.Range("J137").Value = 8500
$ws.Range("L137").Value = 25500
$ws.Range("N137").Value = -35700

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

$ws.Range("H102").Value = 7189
$ws.Range("I102").Value = 7025.125
$ws.Range("K102").Value = 7025.125
$ws.Range("M102").Value = -5403.125

$ws.Range("H126").Value = 8083.039
$ws.Range("I126").Value = 7082.0435
$ws.Range("J126").Value = 8905.286
$ws.Range("K126").Value = 21246.1305
$ws.Range("L126").Value = 26715.858
$ws.Range("M126").Value = -18776.1305
$ws.Range("N126").Value = -31655.858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 52632412
$ws.Range("I16").Value = 880.7222
$ws.Range("J16").Value = 1000000000
$ws.Range("K16").Value = 880.7222
$ws.Range("L16").Value = 1000000000
$ws.Range("M16").Value = -710.7222
$ws.Range("N16").Value = -1000000340

$ws.Range("H40").Value = 6225.3076
$ws.Range("I40").Value = 5471.1763
$ws.Range("K40").Value = 5471.1763
$ws.Range("M40").Value = -5335.1763

$ws.Range("H107").Value = 4099.3335
$ws.Range("I107").Value = 4099.3335
$ws.Range("K107").Value = 4099.3335
$ws.Range("M107").Value = -2179.3335

$ws.Range("H122").Value = 3820.5278
$ws.Range("I122").Value = 3074.4285
$ws.Range("K122").Value = 9223.2855
$ws.Range("M122").Value = -6773.2855

$ws.Range("H132").Value = 10878869
$ws.Range("I132").Value = 17244992
$ws.Range("K132").Value = 51734976
$ws.Range("M132").Value = -51732446

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 889.7857
$ws.Range("I107").Value = 596.875
$ws.Range("J107").Value = 1280.3334
$ws.Range("K107").Value = 1790.625
$ws.Range("L107").Value = 3841.0002
$ws.Range("M107").Value = 129.375
$ws.Range("N107").Value = -7681.0002

$ws.Range("H122").Value = 178429.1
$ws.Range("I122").Value = 310951.47
$ws.Range("K122").Value = 932854.4099999999
$ws.Range("M122").Value = -930404.4099999999

$ws.Range("H132").Value = 4977.025
$ws.Range("I132").Value = 5877.0415
$ws.Range("J132").Value = 3627
$ws.Range("K132").Value = 17631.1245
$ws.Range("L132").Value = 10881
$ws.Range("M132").Value = -15101.1245
$ws.Range("N132").Value = -15941

$ws.Range("H136").Value = 13120221
$ws.Range("I136").Value = 19232776
$ws.Range("J136").Value = 406107.97
$ws.Range("K136").Value = 57698328
$ws.Range("L136").Value = 1218323.91
$ws.Range("M136").Value = -57695778
$ws.Range("N136").Value = -1223423.91

$ws.Range("H141").Value = 90000
$ws.Range("J141").Value = 90000
$ws.Range("L141").Value = 90000
$ws.Range("N141").Value = -100360
